$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix trailing space on A1 ("John " -> "John")
$ws.Range("A1").Value = "John"

# New employee row
$ws.Range("A4").Value = "Aaron"
$ws.Range("B4").Value = "Ho"
$ws.Range("C4").Value = "aaron.ho@email.com"

$ws.Range("D4").Value = 4444444444
$ws.Range("D4").NumberFormat = "0"

$ws.Range("E4").Value = [DateTime]"2024-03-20"

$ws.Range("F4").Value = "Engineer "
$ws.Range("F4").Font.ThemeColor = 1

$ws.Range("G4").Value = 70000
$ws.Range("G4").NumberFormat = "0.00"

$ws.Range("H4").Value = "IT "
$ws.Range("H4").Font.ThemeColor = 1

$ws.Range("I4").Value = "john.doe@email.com "
